$d = $word.ActiveDocument

$replacements = @(
    @("2023-12-18 Monday", "2023-12-19 Tuesday"),
    @("70×90=6300", "14×61=854"),
    @("69×82=5658", "38×79=3002"),
    @("98×92=9016", "82×53=4346"),
    @("97×91=8827", "36×59=2124"),
    @("55×20=1100", "55×36=1980"),
    @("37×90=3330", "52×44=2288"),
    @("27×86=2322", "66×85=5610"),
    @("40×95=3800", "98×41=4018"),
    @("59×38=2242", "43×50=2150"),
    @("45×98=4410", "96×80=7680"),
    @("34×25=850", "71×45=3195"),
    @("20×73=1460", "94×21=1974"),
    @("12×98=1176", "37×37=1369"),
    @("81×27=2187", "55×28=1540"),
    @("99×16=1584", "29×67=1943"),
    @("81×35=2835", "98×53=5194"),
    @("66×54=3564", "44×46=2024"),
    @("25×50=1250", "47×38=1786"),
    @("82×84=6888", "17×86=1462"),
    @("93×28=2604", "80×85=6800"),
    @("70×37=2590", "78×94=7332"),
    @("13×70=910", "44×14=616"),
    @("75×54=4050", "66×63=4158"),
    @("29×82=2378", "98×18=1764"),
    @("25×26=650", "29×27=783")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done."
